$d = $word.ActiveDocument

# Locate the end of the ISP paragraph's text ("... interface genérica.")
$rng = $d.Content
$found = $rng.Find.Execute("interface gen" + [char]0x00E9 + "rica.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate end of ISP paragraph text"
}
$rng.Collapse(0)
$insertAt = $rng.Start

# The trailing bookmark ("_GoBack") currently sits right after this text,
# attached to the ISP paragraph. Remove it here; it gets re-added at the
# new end of document (end of the new last paragraph) below.
try {
    $existing = $d.Bookmarks.Item("_GoBack")
    $existing.Delete()
} catch {
}

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$pPr = "<w:pPr><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr>"

$rPrBold = "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr>"
$rPrPlain = "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr>"

# Paragraph A: empty spacer paragraph
$paraEmpty = "<w:p $w>$pPr</w:p>"

# Paragraph B: heading "5 - DIP - Dependency Inversion Principle: " + definition text
$paraHeading = "<w:p $w>$pPr" +
    "<w:r>$rPrBold<w:t xml:space=`"preserve`">5 " + [char]0x2013 + " DIP " + [char]0x2013 + " </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r>$rPrBold<w:t>Dependency</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r>$rPrBold<w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r>$rPrBold<w:t>Inversion</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r>$rPrBold<w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r>$rPrBold<w:t>Principle</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r>$rPrBold<w:t xml:space=`"preserve`">: </w:t></w:r>" +
    ("<w:r>$rPrPlain<w:t>Princ" + [char]0x00ED + "pio da invers" + [char]0x00E3 + "o de depend" + [char]0x00EA + "ncia. De acordo com ele, m" + [char]0x00F3 + "dulos de alto n" + [char]0x00ED + "vel n" + [char]0x00E3 + "o devem depender de m" + [char]0x00F3 + "dulos de baixo n" + [char]0x00ED + "vel, ambos devem depender da abstra" + [char]0x00E7 + [char]0x00F5 + "es. Al" + [char]0x00E9 + "m disso, abstra" + [char]0x00E7 + [char]0x00F5 + "es n" + [char]0x00E3 + "o devem depender de detalhes, detalhes devem depender de abstra" + [char]0x00E7 + [char]0x00F5 + "es.</w:t></w:r>") +
    "</w:p>"

# Paragraph C: "Não devemos confundir..." paragraph, carries the relocated bookmark
$paraConfundir = "<w:p $w>$pPr" +
    ("<w:r>$rPrPlain<w:tab/><w:t>N" + [char]0x00E3 + "o devemos confundir esse princ" + [char]0x00ED + "pio com a inje" + [char]0x00E7 + [char]0x00E3 + "o de depend" + [char]0x00EA + "ncias. Ambas visam desacoplar o c" + [char]0x00F3 + "digo, por" + [char]0x00E9 + "m s" + [char]0x00E3 + "o coisas diferentes. A inje" + [char]0x00E7 + [char]0x00E3 + "o de depend" + [char]0x00EA + "ncia</w:t></w:r>") +
    ("<w:r>$rPrPlain<w:t xml:space=`"preserve`"> " + [char]0x00E9 + " um padr" + [char]0x00E3 + "o de projeto e a invers" + [char]0x00E3 + "o de depend" + [char]0x00EA + "ncia " + [char]0x00E9 + " um princ" + [char]0x00ED + "pio (conceito).</w:t></w:r>") +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
    "</w:p>"

$xml = $paraEmpty + $paraHeading + $paraConfundir

$insertRange = $d.Range($insertAt, $insertAt)
[void]$insertRange.InsertXML($xml)

$d.Save()
